$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Cells.Item(1, 6).Value = "Relative Importance"

# --- Data rows (row, B, C, D, E, F) ---
$ws.Cells.Item(2, 2).Value = "ln Nphoto"
$ws.Cells.Item(2, 3).Value = 0.51999005635377105
$ws.Cells.Item(2, 4).Value = [double]"9.1451990981480297E-2"
$ws.Cells.Item(2, 5).Value = [double]"1.00451004079254E-8"
$ws.Cells.Item(2, 6).Value = 23.248516437693599

$ws.Cells.Item(3, 2).Value = "ln Nstructure"
$ws.Cells.Item(3, 3).Value = 0.95191716087175204
$ws.Cells.Item(3, 4).Value = [double]"8.5038995856662695E-3"
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 37.101807816901101

$ws.Cells.Item(4, 2).Value = "Soil N"
$ws.Cells.Item(4, 3).Value = "NA"
$ws.Cells.Item(4, 4).Value = "NA"
$ws.Cells.Item(4, 5).Value = [double]"2.81177763912933E-58"
$ws.Cells.Item(4, 6).Value = 5.2283856841646204

$ws.Cells.Item(5, 2).Value = "Soil P"
$ws.Cells.Item(5, 3).Value = "NA"
$ws.Cells.Item(5, 4).Value = "NA"
$ws.Cells.Item(5, 5).Value = 0.71885048633022097
$ws.Cells.Item(5, 6).Value = 3.8214794973311301

$ws.Cells.Item(6, 2).Value = "Soil K+¬µ"
$ws.Cells.Item(6, 3).Value = "NA"
$ws.Cells.Item(6, 4).Value = "NA"
$ws.Cells.Item(6, 5).Value = 0.50568501053739001
$ws.Cells.Item(6, 6).Value = 3.5335695250771799

$ws.Cells.Item(7, 2).Value = "N fixer"
$ws.Cells.Item(7, 3).Value = "NA"
$ws.Cells.Item(7, 4).Value = "NA"
$ws.Cells.Item(7, 5).Value = [double]"8.1035966827340394E-17"
$ws.Cells.Item(7, 6).Value = 4.4190929997844997

$ws.Cells.Item(8, 2).Value = "C3/C4"
$ws.Cells.Item(8, 3).Value = "NA"
$ws.Cells.Item(8, 4).Value = "NA"
$ws.Cells.Item(8, 5).Value = [double]"3.0384325381465902E-21"
$ws.Cells.Item(8, 6).Value = 8.9262544586997894

$ws.Cells.Item(9, 2).Value = "Soil N x Soil P"
$ws.Cells.Item(9, 3).Value = "NA"
$ws.Cells.Item(9, 4).Value = "NA"
$ws.Cells.Item(9, 5).Value = [double]"3.2949809096825601E-3"
$ws.Cells.Item(9, 6).Value = 1.2763755294413299

$ws.Cells.Item(10, 2).Value = "Soil N x Soil P"
$ws.Cells.Item(10, 3).Value = "NA"
$ws.Cells.Item(10, 4).Value = "NA"
$ws.Cells.Item(10, 5).Value = 0.59579950859943098
$ws.Cells.Item(10, 6).Value = 1.19035243608916

$ws.Cells.Item(11, 2).Value = "Soil P x Soil K"
$ws.Cells.Item(11, 3).Value = "NA"
$ws.Cells.Item(11, 4).Value = "NA"
$ws.Cells.Item(11, 5).Value = 0.80298402508854305
$ws.Cells.Item(11, 6).Value = 1.11111870647406

$ws.Cells.Item(12, 2).Value = "Soil N x Soil P x Soil K"
$ws.Cells.Item(12, 3).Value = "NA"
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(12, 5).Value = 0.98116941100215505
$ws.Cells.Item(12, 6).Value = 0.60164918721258598

# --- Apply the new "0.00" number format to column F rows 2-12 (creates a new cellXf, numFmtId 2) ---
$ws.Range("F2:F12").NumberFormat = "0.00"

# --- Update selection to F2 ---
$ws.Range("F2").Select()
